$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.803.54"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "1.583.43"
$ws.Range("E3").Value = "  -2.22%  "
$ws.Range("E4").Value = "  -0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.19"
$ws.Range("E5").Value = "  -1.73%  "
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("E7").Value = "  -3.63%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.247"
$ws.Range("E8").Value = "  -0.72%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0618"
$ws.Range("E9").Value = "  +0.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.19"
$ws.Range("E10").Value = "  -1.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0787"
$ws.Range("E11").Value = "  -0.38%  "
$ws.Range("D12").Value = "1.802.16"
$ws.Range("E12").Value = "  -2.27%  "
$ws.Range("D13").Value = "1.592.38"
$ws.Range("E13").Value = "  -1.57%  "
$ws.Range("E14").Value = "  -2.80%  "
$ws.Range("E15").Value = "  -2.21%  "
$ws.Range("D16").Value = "25.786.84"
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.25"
$ws.Range("E17").Value = "  -2.04%  "
$ws.Range("D18").Value = "0.0₃0722"
$ws.Range("E18").Value = "  -1.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.00"
$ws.Range("E19").Value = "  -0.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "191.76"
$ws.Range("E20").Value = "  +0.25%  "
$ws.Range("E21").Value = "  -1.22%  "
$ws.Range("E22").Value = "  -0.78%  "
$ws.Range("E23").Value = "  -1.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.132"
$ws.Range("E24").Value = "  -2.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.25"
$ws.Range("E25").Value = "  -1.71%  "
$ws.Range("E26").Value = "  -0.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.70"
$ws.Range("E27").Value = "  -1.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.15"
$ws.Range("E28").Value = "  -0.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.47"
$ws.Range("E29").Value = "  -2.68%  "
$ws.Range("E30").Value = "  -5.89%  "
$ws.Range("E31").Value = "  -0.92%  "
$ws.Range("E32").Value = "  -0.46%  "
$ws.Range("E34").Value = "  +0.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.32"
$ws.Range("E35").Value = "  -3.98%  "
$ws.Range("D36").Value = "1.104.49"
$ws.Range("E36").Value = "  -1.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.53%  "
$ws.Range("E38").Value = "  -1.27%  "
$ws.Range("E39").Value = "  -2.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0150"
$ws.Range("E40").Value = "  -1.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.783"
$ws.Range("E41").Value = "  -6.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.809"
$ws.Range("E42").Value = "  +8.18%  "
$ws.Range("E43").Value = "  +2.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "93.28"
$ws.Range("E44").Value = "  -5.07%  "
$ws.Range("D45").Value = "1.716.58"
$ws.Range("E45").Value = "  -2.21%  "
$ws.Range("D46").Value = "0.0₆0112"
$ws.Range("E46").Value = "  -0.97%  "
$ws.Range("E47").Value = "  -0.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.28"
$ws.Range("E48").Value = "  -1.34%  "
$ws.Range("E49").Value = "  -1.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.405"
$ws.Range("E50").Value = "  -1.33%  "
$ws.Range("E51").Value = "  -0.52%  "
